$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.005.22"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "1.640.65"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  -0.82%  "

$ws.Range("D5").Value = "215.72"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "0.5153"
$ws.Range("E6").Value = "  +1.39%  "

$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("D8").Value = "0.2587"
$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("D9").Value = "0.06377"
$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").Value = "19.85"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").Value = "0.07771"
$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("D12").Value = "4.292"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").Value = "1.638.99"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").Value = "0.5481"
$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("D15").Value = "0.0₅7784"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").Value = "64.57"
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").Value = "26.032.21"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").Value = "199.57"
$ws.Range("E19").Value = "  +1.35%  "

$ws.Range("D20").Value = "4.469"
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("D21").Value = "9.986"
$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("D22").Value = "6.111"
$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("D24").Value = "1.902"
$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("D25").Value = "142.33"
$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").Value = "0.1234"
$ws.Range("E26").Value = "  +7.62%  "

$ws.Range("D27").Value = "6.877"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").Value = "15.66"
$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("D29").Value = "1.244"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").Value = "0.04868"
$ws.Range("E30").Value = "  -3.27%  "

$ws.Range("D31").Value = "3.312"
$ws.Range("E31").Value = "  +1.09%  "

$ws.Range("D32").Value = "3.232"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").Value = "1.543"
$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").Value = "2.377"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").Value = "0.9191"
$ws.Range("E35").Value = "  +2.74%  "

$ws.Range("D36").Value = "0.5599"
$ws.Range("E36").Value = "  +0.91%  "

$ws.Range("D37").Value = "2.570"
$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("D38").Value = "1.123.65"
$ws.Range("E38").Value = "  -0.85%  "

$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("D40").Value = "1.001"

$ws.Range("D41").Value = "2.526"
$ws.Range("E41").Value = "  -1.32%  "

$ws.Range("D42").Value = "5.584"
$ws.Range("E42").Value = "  -1.51%  "

$ws.Range("D43").Value = "0.8095"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("D44").Value = "99.74"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").Value = "1.780.44"
$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("E47").Value = "  -0.31%  "

$ws.Range("D48").Value = "55.38"
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("D49").Value = "1.007"
$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("D50").Value = "0.05221"
$ws.Range("E50").Value = "  +2.40%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.512"
$ws.Range("E51").Value = "  +0.97%  "
